$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.547.38"
$ws.Range("E2").Value = "  +6.47%  "
$ws.Range("D3").Value = "2.419.65"
$ws.Range("E3").Value = "  +6.17%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'473.23"
$ws.Range("E5").Value = "  +10.58%  "
$ws.Range("D6").Value = "'137.30"
$ws.Range("E6").Value = "  +15.65%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'0.499"
$ws.Range("E8").Value = "  +9.55%  "
$ws.Range("D9").Value = "2.444.83"
$ws.Range("E9").Value = "  +7.36%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'5.49"
$ws.Range("E10").Value = "  +9.06%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0949"
$ws.Range("E11").Value = "  +8.63%  "
$ws.Range("D12").Value = "'0.321"
$ws.Range("E12").Value = "  +7.63%  "
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").Value = "2.843.29"
$ws.Range("E14").Value = "  +7.10%  "
$ws.Range("D15").Value = "54.724.40"
$ws.Range("E15").Value = "  +7.04%  "
$ws.Range("D16").Value = "'20.17"
$ws.Range("E16").Value = "  +9.17%  "
$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = "  +13.68%  "
$ws.Range("D18").Value = "2.434.30"
$ws.Range("E18").Value = "  +7.00%  "
$ws.Range("D19").Value = "'4.29"
$ws.Range("E19").Value = "  +7.68%  "
$ws.Range("D20").Value = "'9.77"
$ws.Range("E20").Value = "  +12.83%  "
$ws.Range("D21").Value = "'309.97"
$ws.Range("E21").Value = "  +6.43%  "
$ws.Range("D22").Value = "'0.996"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "'5.63"
$ws.Range("E23").Value = "  +10.69%  "
$ws.Range("D24").Value = "'56.59"
$ws.Range("E24").Value = "  +8.14%  "
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "'0.397"
$ws.Range("E26").Value = "  +8.52%  "
$ws.Range("D27").Value = "'0.159"
$ws.Range("E27").Value = "  +20.39%  "
$ws.Range("D28").Value = "2.535.19"
$ws.Range("E28").Value = "  +9.38%  "
$ws.Range("D29").Value = "'7.20"
$ws.Range("E29").Value = "  +8.17%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0754"
$ws.Range("E31").Value = "  +16.15%  "
$ws.Range("D32").Value = "'147.89"
$ws.Range("E32").Value = "  +3.78%  "
$ws.Range("D33").Value = "'17.95"
$ws.Range("E33").Value = "  +8.02%  "
$ws.Range("D34").Value = "'1.45"
$ws.Range("E34").Value = "  +11.67%  "
$ws.Range("D35").Value = "'5.07"
$ws.Range("E35").Value = "  +7.48%  "
$ws.Range("D36").Value = "'1.11"
$ws.Range("E36").Value = "  +14.01%  "
$ws.Range("D37").Value = "'3.52"
$ws.Range("E37").Value = "  +8.93%  "
$ws.Range("D38").Value = "'0.821"
$ws.Range("E38").Value = "  +11.33%  "
$ws.Range("D39").Value = "'33.47"
$ws.Range("E39").Value = "  +5.30%  "
$ws.Range("D40").Value = "'0.995"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").Value = "'3.39"
$ws.Range("E41").Value = "  +8.87%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.596"
$ws.Range("E42").Value = "  +8.38%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0539"
$ws.Range("E43").Value = "  +9.34%  "
$ws.Range("D44").Value = "'1.26"
$ws.Range("E44").Value = "  +12.30%  "
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'252.60"
$ws.Range("E46").Value = "  +33.67%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'4.57"
$ws.Range("E47").Value = "  +19.34%  "
$ws.Range("D48").Value = "'0.0882"
$ws.Range("E48").Value = "  +10.75%  "
$ws.Range("D49").Value = "'0.0219"
$ws.Range("E49").Value = "  +9.46%  "
$ws.Range("D50").Value = "1.886.59"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "'16.84"
$ws.Range("E51").Value = "  +8.57%  "
